# More tests for Persons obj.
# - Fix typo in shared string used by B3: "Mantainence" -> "Maintainence"
# - Rename the worksheet tab: "Blad1" -> "test_departments"
# - Move the active selection to B3
# - Touch page setup (paper size / orientation) so a <pageSetup> element is emitted

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Correct the typo in the "Maint" synonym cell.
$ws.Range("B3").Value = "Maintainence"

# Rename the sheet to match the data file name.
$ws.Name = "test_departments"

# Update the saved selection/active cell.
$ws.Range("B3").Select() | Out-Null

# Set print/page setup (adds a <pageSetup> element on save).
$ws.PageSetup.PaperSize = 9
$ws.PageSetup.Orientation = 1
